# Cotações atualizadas - 2025-09-27
# Append the new daily quote row (row 23) below the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data (serial date 45927 = 2025-09-27) - match the number format used by
# the existing date column (copy from the row above it).
$ws.Range("A23").Value = 45927
$ws.Range("A23").NumberFormat = $ws.Range("A22").NumberFormat

# Fund quotes for the day, stored as text (comma decimal separator),
# matching the existing rows' formatting.
$ws.Range("B23").Value = "21,0192"
$ws.Range("C23").Value = "15,0785"
$ws.Range("D23").Value = "14,8412"
$ws.Range("E23").Value = "14,8412"
